$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13, shifting existing rows 13:73 down to 15:75
$ws.Rows("13:14").Insert()

# Row 13: new "Primera" record (date 2021-12-20 serial 44550)
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44550
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100114007
$ws.Range("G13").Value = "Jengibre"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 790
$ws.Range("K13").Value = 11000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 11494
$ws.Range("N13").Value = "$/caja 13 kilos"
$ws.Range("O13").Value = "Perú"
$ws.Range("P13").Value = 884
$ws.Range("Q13").Value = 13
$ws.Range("R13").Value = "Hortaliza"

# Row 14: new "Segunda" record (same date)
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = 44550
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100114007
$ws.Range("G14").Value = "Jengibre"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 430
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9500
$ws.Range("N14").Value = "$/caja 13 kilos"
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 731
$ws.Range("Q14").Value = 13
$ws.Range("R14").Value = "Hortaliza"
